# Applies the "Updated symbol list" crypto data refresh to Sheet1.
# Each entry updates a single cell while preserving it as text
# (values like "320.00" or "-1.57%" must stay text, matching the
# original inline-string cells, instead of being auto-converted by
# Excel into numbers/percentages).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "320.00" },
    @{ Cell = "E2"; Value = "-1.57%" },
    @{ Cell = "D3"; Value = "39.32" },
    @{ Cell = "E3"; Value = "-1.36%" },
    @{ Cell = "D4"; Value = "5.884" },
    @{ Cell = "E4"; Value = "12.49%" },
    @{ Cell = "D5"; Value = "0.08010" },
    @{ Cell = "E5"; Value = "-1.02%" },
    @{ Cell = "B6"; Value = "KuCoinToken" },
    @{ Cell = "C6"; Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs" },
    @{ Cell = "D6"; Value = "8.634" },
    @{ Cell = "E6"; Value = "0.07%" },
    @{ Cell = "B7"; Value = "FTXToken" },
    @{ Cell = "C7"; Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt" },
    @{ Cell = "D7"; Value = "1.894" },
    @{ Cell = "E7"; Value = "-1.62%" },
    @{ Cell = "B8"; Value = "MXToken" },
    @{ Cell = "C8"; Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx" },
    @{ Cell = "D8"; Value = "0.9352" },
    @{ Cell = "E8"; Value = "0.15%" },
    @{ Cell = "B9"; Value = "LiechtensteinCryptoassetsExchange" },
    @{ Cell = "C9"; Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx" },
    @{ Cell = "D9"; Value = "0.1250" },
    @{ Cell = "E9"; Value = "-3.87%" },
    @{ Cell = "B10"; Value = "WazirX" },
    @{ Cell = "C10"; Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx" },
    @{ Cell = "D10"; Value = "0.1948" },
    @{ Cell = "E10"; Value = "-0.35%" },
    @{ Cell = "B11"; Value = "MCDex" },
    @{ Cell = "C11"; Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb" },
    @{ Cell = "D11"; Value = "8.751" },
    @{ Cell = "E11"; Value = "30.68%" },
    @{ Cell = "B12"; Value = "MandalaExchangeToken" },
    @{ Cell = "C12"; Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx" },
    @{ Cell = "D12"; Value = "0.09117" },
    @{ Cell = "E12"; Value = "-0.42%" },
    @{ Cell = "B13"; Value = "BitrueCoin" },
    @{ Cell = "C13"; Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr" },
    @{ Cell = "D13"; Value = "0.03510" },
    @{ Cell = "E13"; Value = "3.11%" },
    @{ Cell = "B14"; Value = "BitMartToken" },
    @{ Cell = "C14"; Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx" },
    @{ Cell = "D14"; Value = "0.09571" },
    @{ Cell = "E14"; Value = "0.31%" },
    @{ Cell = "B15"; Value = "BitForexToken" },
    @{ Cell = "C15"; Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf" },
    @{ Cell = "D15"; Value = "0.001282" },
    @{ Cell = "E15"; Value = "-7.75%" },
    @{ Cell = "B16"; Value = "TigerCash" },
    @{ Cell = "C16"; Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch" },
    @{ Cell = "D16"; Value = "0.006103" },
    @{ Cell = "E16"; Value = "-5.11%" },
    @{ Cell = "B17"; Value = "LEO" },
    @{ Cell = "C17"; Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo" },
    @{ Cell = "D17"; Value = "3.355" },
    @{ Cell = "E17"; Value = "-0.07%" },
    @{ Cell = "B18"; Value = "GateToken" },
    @{ Cell = "C18"; Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt" },
    @{ Cell = "D18"; Value = "4.572" },
    @{ Cell = "E18"; Value = "0.88%" },
    @{ Cell = "D19"; Value = "2.951" },
    @{ Cell = "E19"; Value = "-0.26%" },
    @{ Cell = "E20"; Value = "0.01%" },
    @{ Cell = "D21"; Value = "0.1430" },
    @{ Cell = "E21"; Value = "7.67%" },
    @{ Cell = "D22"; Value = "0.2412" },
    @{ Cell = "E22"; Value = "4.37%" },
    @{ Cell = "D23"; Value = "0.04465" },
    @{ Cell = "E23"; Value = "0.76%" },
    @{ Cell = "D24"; Value = "0.001264" },
    @{ Cell = "E24"; Value = "3.41%" },
    @{ Cell = "D25"; Value = "0.004417" },
    @{ Cell = "E25"; Value = "1.39%" },
    @{ Cell = "D26"; Value = "0.0001142" },
    @{ Cell = "E26"; Value = "-11.44%" },
    @{ Cell = "E27"; Value = "0.11%" },
    @{ Cell = "D39"; Value = "0.02397" },
    @{ Cell = "E39"; Value = "-2.68%" },
    @{ Cell = "D40"; Value = "0.05174" },
    @{ Cell = "E40"; Value = "-0.88%" },
    @{ Cell = "D41"; Value = "0.007431" },
    @{ Cell = "E41"; Value = "-3.27%" },
    @{ Cell = "B42"; Value = "Dexo" },
    @{ Cell = "C42"; Value = "https://coinranking.com/coin/QkL_pl546+dexo-dexo" },
    @{ Cell = "D42"; Value = "0.009155" },
    @{ Cell = "E42"; Value = "6.09%" },
    @{ Cell = "B43"; Value = "BKEXToken" },
    @{ Cell = "C43"; Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk" },
    @{ Cell = "D43"; Value = "0.1402" },
    @{ Cell = "E43"; Value = "-2.11%" },
    @{ Cell = "D44"; Value = "0.002124" },
    @{ Cell = "E44"; Value = "0.69%" },
    @{ Cell = "D45"; Value = "0.01116" },
    @{ Cell = "E45"; Value = "37.33%" },
    @{ Cell = "D46"; Value = "0.00006743" },
    @{ Cell = "E46"; Value = "2.03%" },
    @{ Cell = "D47"; Value = "0.00000000751" },
    @{ Cell = "E47"; Value = "0.21%" },
    @{ Cell = "D48"; Value = "0.003010" },
    @{ Cell = "E48"; Value = "5.60%" },
    @{ Cell = "D50"; Value = "0.00002104" },
    @{ Cell = "E50"; Value = "0.21%" },
    @{ Cell = "D51"; Value = "0.0002004" },
    @{ Cell = "E51"; Value = "0.21%" }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    # Force text format so numeric-looking / percent-looking strings
    # are not reinterpreted as numbers by Excel.
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    # Drop the explicit text-format styling again so the cell keeps
    # using the default (unstyled) look, as in the original file.
    $rng.ClearFormats()
}
